$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.191.48'
$ws.Range('E2').Value = '  +0.61%  '

$ws.Range('D3').Value = '2.569.97'
$ws.Range('E3').Value = '  +1.09%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = "'585.28"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.30%  '

$ws.Range('D6').Value = "'147.33"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.45%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  +3.06%  '

$ws.Range('E9').Value = '  +3.65%  '

$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('D12').Value = "'0.356"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.27%  '

$ws.Range('D13').Value = "'27.38"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.69%  '

$ws.Range('D14').Value = '3.030.65'
$ws.Range('E14').Value = '  +1.08%  '

$ws.Range('D15').Value = '63.158.92'
$ws.Range('E15').Value = '  +0.52%  '

$ws.Range('E16').Value = '  +3.94%  '

$ws.Range('D17').Value = '2.569.25'
$ws.Range('E17').Value = '  +0.81%  '

$ws.Range('E18').Value = '  -0.71%  '

$ws.Range('D19').Value = "'344.25"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.55%  '

$ws.Range('D20').Value = "'4.42"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.12%  '

$ws.Range('D21').Value = "'6.89"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.29%  '

$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('B23').Value = 'LEO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D23').Value = "'5.54"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.70%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'66.91"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.41%  '

$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '2.697.19'
$ws.Range('E25').Value = '  +0.95%  '

$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = "'0.171"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.35%  '

$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').Value = "'1.63"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.53%  '

$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').Value = "'8.14"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +12.05%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = "'8.50"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.27%  '

$ws.Range('B30').Value = 'SuiNetwork'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D30').Value = "'1.49"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.46%  '

$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = "'0.999"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.02%  '

$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = "'1.99"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.04%  '

$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '0.0₃0827'
$ws.Range('E33').Value = '  +2.27%  '

$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = "'467.69"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.30%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = "'1.64"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.21%  '

$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = "'176.13"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.54%  '

$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = "'0.409"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.67%  '

$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = "'19.23"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.38%  '

$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = "'4.55"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.37%  '

$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').Value = "'0.999"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = "'1.75"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.48%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '

$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = "'151.88"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.81%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = "'3.81"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.07%  '

$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = "'20.99"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.19%  '

$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = "'0.0547"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.67%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = "'0.613"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.21%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = "'0.0979"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.36%  '

$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = "'0.0239"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.05%  '

$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = "'1.75"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.37%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = "'11.38"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.09%  '
